# Update ShipmentTracking (column P) and ActualRate (column Q, row 24) values
# with new FedEx tracking numbers / rate, per "final changes 17th march 2022".
#
# Values are written via a temporary formula (="...") and then converted to a
# plain static value with PasteSpecial (values only). This keeps the cell as a
# text value (matching the original shared-string text cells, since these
# tracking numbers / "$" amounts would otherwise be auto-converted to numbers
# by Excel) without leaving a formula behind or mutating the cell's style.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("P2").Formula = '="320018063820"'
$ws.Range("P2").Copy()
$ws.Range("P2").PasteSpecial(-4163)

$ws.Range("P3").Formula = '="320018063831"'
$ws.Range("P3").Copy()
$ws.Range("P3").PasteSpecial(-4163)

$ws.Range("P4").Formula = '="320018063864"'
$ws.Range("P4").Copy()
$ws.Range("P4").PasteSpecial(-4163)

$ws.Range("P5").Formula = '="320018063886"'
$ws.Range("P5").Copy()
$ws.Range("P5").PasteSpecial(-4163)

$ws.Range("P6").Formula = '="320018063923"'
$ws.Range("P6").Copy()
$ws.Range("P6").PasteSpecial(-4163)

$ws.Range("P7").Formula = '="320018063945"'
$ws.Range("P7").Copy()
$ws.Range("P7").PasteSpecial(-4163)

$ws.Range("P8").Formula = '="320018063978"'
$ws.Range("P8").Copy()
$ws.Range("P8").PasteSpecial(-4163)

$ws.Range("P9").Formula = '="320018063990"'
$ws.Range("P9").Copy()
$ws.Range("P9").PasteSpecial(-4163)

$ws.Range("P10").Formula = '="320018064025"'
$ws.Range("P10").Copy()
$ws.Range("P10").PasteSpecial(-4163)

$ws.Range("P11").Formula = '="320018064047"'
$ws.Range("P11").Copy()
$ws.Range("P11").PasteSpecial(-4163)

$ws.Range("P12").Formula = '="320018064080"'
$ws.Range("P12").Copy()
$ws.Range("P12").PasteSpecial(-4163)

$ws.Range("P13").Formula = '="320018064106"'
$ws.Range("P13").Copy()
$ws.Range("P13").PasteSpecial(-4163)

$ws.Range("P14").Formula = '="320018064139"'
$ws.Range("P14").Copy()
$ws.Range("P14").PasteSpecial(-4163)

$ws.Range("P15").Formula = '="320018064150"'
$ws.Range("P15").Copy()
$ws.Range("P15").PasteSpecial(-4163)

$ws.Range("P16").Formula = '="320018064183"'
$ws.Range("P16").Copy()
$ws.Range("P16").PasteSpecial(-4163)

$ws.Range("P17").Formula = '="320018064209"'
$ws.Range("P17").Copy()
$ws.Range("P17").PasteSpecial(-4163)

$ws.Range("P18").Formula = '="320018064242"'
$ws.Range("P18").Copy()
$ws.Range("P18").PasteSpecial(-4163)

$ws.Range("P19").Formula = '="320018064264"'
$ws.Range("P19").Copy()
$ws.Range("P19").PasteSpecial(-4163)

$ws.Range("P20").Formula = '="320018064297"'
$ws.Range("P20").Copy()
$ws.Range("P20").PasteSpecial(-4163)

$ws.Range("P21").Formula = '="320018064312"'
$ws.Range("P21").Copy()
$ws.Range("P21").PasteSpecial(-4163)

$ws.Range("P22").Formula = '="320018064345"'
$ws.Range("P22").Copy()
$ws.Range("P22").PasteSpecial(-4163)

$ws.Range("P23").Formula = '="320018064356"'
$ws.Range("P23").Copy()
$ws.Range("P23").PasteSpecial(-4163)

$ws.Range("P24").Formula = '="320018064367"'
$ws.Range("P24").Copy()
$ws.Range("P24").PasteSpecial(-4163)

$ws.Range("Q24").Formula = '="$248.51"'
$ws.Range("Q24").Copy()
$ws.Range("Q24").PasteSpecial(-4163)

$ws.Range("P25").Formula = '="320018064378"'
$ws.Range("P25").Copy()
$ws.Range("P25").PasteSpecial(-4163)
